# Apply updated cryptocurrency price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.782.32"
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").Value = "2.291.33"
$ws.Range("E3").Value = "  -4.38%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'533.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.49%  "
$ws.Range("D6").Value = "'130.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("D8").Value = "'0.576"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "2.284.19"
$ws.Range("E9").Value = "  -4.58%  "
$ws.Range("D10").Value = "'0.0995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.54%  "
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").Value = "'23.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.81%  "
$ws.Range("D15").Value = "2.695.02"
$ws.Range("E15").Value = "  -4.59%  "
$ws.Range("D16").Value = "57.777.13"
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("E17").Value = "  -4.43%  "
$ws.Range("D18").Value = "2.314.09"
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("D19").Value = "'10.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.53%  "
$ws.Range("E20").Value = "  -6.32%  "
$ws.Range("D21").Value = "'313.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("E22").Value = "  -5.54%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'62.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("D25").Value = "'0.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.77%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'7.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.22%  "
$ws.Range("E28").Value = "  -6.34%  "
$ws.Range("D29").Value = "'170.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  -6.13%  "
$ws.Range("D31").Value = "0.0₃0718"
$ws.Range("E31").Value = "  -6.46%  "
$ws.Range("E32").Value = "  -6.06%  "
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("E34").Value = "  -5.89%  "
$ws.Range("D36").Value = "'17.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  -7.47%  "
$ws.Range("D39").Value = "'3.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.67%  "
$ws.Range("D40").Value = "'38.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").Value = "'1.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.64%  "
$ws.Range("D42").Value = "'141.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.93%  "
$ws.Range("D43").Value = "'288.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.30%  "
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("E46").Value = "  -3.53%  "
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("D48").Value = "'18.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.78%  "
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "0.0₆0202"
$ws.Range("E51").Value = "  +81.73%  "
